$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 38 <- original row 39 data (columns B:AD); column A (id) stays fixed
$ws.Range('B38').Value2 = 6781354
$ws.Range('C38').Value2 = 'Costa Rica Primera Division'
$ws.Range('D38').Value2 = 45171.75
$ws.Range('E38').Value2 = 'Puntarenas'
$ws.Range('F38').Value2 = 'AD San Carlos'
$ws.Range('G38').Value2 = 1
$ws.Range('H38').Value2 = 0
$ws.Range('I38').Value2 = 0
$ws.Range('J38').Value2 = 0
$ws.Range('K38').Value2 = 'H'
$ws.Range('L38').Value2 = 2.4
$ws.Range('M38').Value2 = 3.2
$ws.Range('N38').Value2 = 2.8
$ws.Range('O38').Value2 = 2.3
$ws.Range('P38').Value2 = 3.2
$ws.Range('Q38').Value2 = 3
$ws.Range('R38').Value2 = -0.25
$ws.Range('S38').Value2 = 2
$ws.Range('T38').Value2 = 1.8
$ws.Range('U38').Value2 = 2.25
$ws.Range('V38').Value2 = 1.9
$ws.Range('W38').Value2 = 1.9
$ws.Range('X38').Value2 = 1.3
$ws.Range('Y38').Value2 = -1
$ws.Range('Z38').Value2 = -1
$ws.Range('AA38').Value2 = 1
$ws.Range('AB38').Value2 = -1
$ws.Range('AC38').Value2 = -1
$ws.Range('AD38').Value2 = 0.8999999999999999

# Row 39 <- original row 38 data (columns B:AD); column A (id) stays fixed
$ws.Range('B39').Value2 = 6782522
$ws.Range('C39').Value2 = 'Costa Rica Primera Division'
$ws.Range('D39').Value2 = 45171.75
$ws.Range('E39').Value2 = 'Municipal Perez Zeledon'
$ws.Range('F39').Value2 = 'Sporting San Jose'
$ws.Range('G39').Value2 = 1
$ws.Range('H39').Value2 = 2
$ws.Range('I39').Value2 = 0
$ws.Range('J39').Value2 = 1
$ws.Range('K39').Value2 = 'A'
$ws.Range('L39').Value2 = 2.5
$ws.Range('M39').Value2 = 3.5
$ws.Range('N39').Value2 = 2.5
$ws.Range('O39').Value2 = 2.2
$ws.Range('P39').Value2 = 3.5
$ws.Range('Q39').Value2 = 2.9
$ws.Range('R39').Value2 = -0.25
$ws.Range('S39').Value2 = 1.9
$ws.Range('T39').Value2 = 1.9
$ws.Range('U39').Value2 = 2.5
$ws.Range('V39').Value2 = 1.9
$ws.Range('W39').Value2 = 1.9
$ws.Range('X39').Value2 = -1
$ws.Range('Y39').Value2 = -1
$ws.Range('Z39').Value2 = 1.9
$ws.Range('AA39').Value2 = -1
$ws.Range('AB39').Value2 = 0.8999999999999999
$ws.Range('AC39').Value2 = 0.8999999999999999
$ws.Range('AD39').Value2 = -1

# Row 95 <- original row 96 data (columns B:AD); column A (id) stays fixed
$ws.Range('B95').Value2 = 6782567
$ws.Range('C95').Value2 = 'Costa Rica Primera Division'
$ws.Range('D95').Value2 = 45221.79166666666
$ws.Range('E95').Value2 = 'AD Grecia'
$ws.Range('F95').Value2 = 'Municipal Liberia'
$ws.Range('G95').Value2 = 2
$ws.Range('H95').Value2 = 3
$ws.Range('I95').Value2 = 2
$ws.Range('J95').Value2 = 2
$ws.Range('K95').Value2 = 'A'
$ws.Range('L95').Value2 = 2.875
$ws.Range('M95').Value2 = 3.5
$ws.Range('N95').Value2 = 2.15
$ws.Range('O95').Value2 = 2.3
$ws.Range('P95').Value2 = 3.5
$ws.Range('Q95').Value2 = 2.6
$ws.Range('R95').Value2 = 0
$ws.Range('S95').Value2 = 1.8
$ws.Range('T95').Value2 = 2
$ws.Range('U95').Value2 = 2.75
$ws.Range('V95').Value2 = 1.8
$ws.Range('W95').Value2 = 2
$ws.Range('X95').Value2 = -1
$ws.Range('Y95').Value2 = -1
$ws.Range('Z95').Value2 = 1.6
$ws.Range('AA95').Value2 = -1
$ws.Range('AB95').Value2 = 1
$ws.Range('AC95').Value2 = 0.8
$ws.Range('AD95').Value2 = -1

# Row 96 <- original row 95 data (columns B:AD); column A (id) stays fixed
$ws.Range('B96').Value2 = 6782565
$ws.Range('C96').Value2 = 'Costa Rica Primera Division'
$ws.Range('D96').Value2 = 45221.79166666666
$ws.Range('E96').Value2 = 'Santos de Gupiles'
$ws.Range('F96').Value2 = 'Municipal Perez Zeledon'
$ws.Range('G96').Value2 = 2
$ws.Range('H96').Value2 = 0
$ws.Range('I96').Value2 = 1
$ws.Range('J96').Value2 = 0
$ws.Range('K96').Value2 = 'H'
$ws.Range('L96').Value2 = 1.833
$ws.Range('M96').Value2 = 3.4
$ws.Range('N96').Value2 = 3.75
$ws.Range('O96').Value2 = 1.833
$ws.Range('P96').Value2 = 3.5
$ws.Range('Q96').Value2 = 3.75
$ws.Range('R96').Value2 = -0.5
$ws.Range('S96').Value2 = 1.875
$ws.Range('T96').Value2 = 1.925
$ws.Range('U96').Value2 = 2.75
$ws.Range('V96').Value2 = 2
$ws.Range('W96').Value2 = 1.8
$ws.Range('X96').Value2 = 0.833
$ws.Range('Y96').Value2 = -1
$ws.Range('Z96').Value2 = -1
$ws.Range('AA96').Value2 = 0.875
$ws.Range('AB96').Value2 = -1
$ws.Range('AC96').Value2 = -1
$ws.Range('AD96').Value2 = 0.8

# Row 129 <- original row 130 data (columns B:AD); column A (id) stays fixed
$ws.Range('B129').Value2 = 6782598
$ws.Range('C129').Value2 = 'Costa Rica Primera Division'
$ws.Range('D129').Value2 = 45255.95833333334
$ws.Range('E129').Value2 = 'Municipal Perez Zeledon'
$ws.Range('F129').Value2 = 'Cartagines'
$ws.Range('G129').Value2 = 1
$ws.Range('H129').Value2 = 0
$ws.Range('I129').Value2 = 1
$ws.Range('J129').Value2 = 0
$ws.Range('K129').Value2 = 'H'
$ws.Range('L129').Value2 = 4.5
$ws.Range('M129').Value2 = 3.75
$ws.Range('N129').Value2 = 1.615
$ws.Range('O129').Value2 = 3.4
$ws.Range('P129').Value2 = 3.4
$ws.Range('Q129').Value2 = 1.85
$ws.Range('R129').Value2 = 0.5
$ws.Range('S129').Value2 = 1.8
$ws.Range('T129').Value2 = 2
$ws.Range('U129').Value2 = 2.75
$ws.Range('V129').Value2 = 1.9
$ws.Range('W129').Value2 = 1.9
$ws.Range('X129').Value2 = 2.4
$ws.Range('Y129').Value2 = -1
$ws.Range('Z129').Value2 = -1
$ws.Range('AA129').Value2 = 0.8
$ws.Range('AB129').Value2 = -1
$ws.Range('AC129').Value2 = -1
$ws.Range('AD129').Value2 = 0.8999999999999999

# Row 130 <- original row 131 data (columns B:AD); column A (id) stays fixed
$ws.Range('B130').Value2 = 6782595
$ws.Range('C130').Value2 = 'Costa Rica Primera Division'
$ws.Range('D130').Value2 = 45255.95833333334
$ws.Range('E130').Value2 = 'Herediano'
$ws.Range('F130').Value2 = 'Sporting San Jose'
$ws.Range('G130').Value2 = 3
$ws.Range('H130').Value2 = 0
$ws.Range('I130').Value2 = 0
$ws.Range('J130').Value2 = 0
$ws.Range('K130').Value2 = 'H'
$ws.Range('L130').Value2 = 1.4
$ws.Range('M130').Value2 = 4.75
$ws.Range('N130').Value2 = 7
$ws.Range('O130').Value2 = 1.363
$ws.Range('P130').Value2 = 4.75
$ws.Range('Q130').Value2 = 8.5
$ws.Range('R130').Value2 = -1.25
$ws.Range('S130').Value2 = 1.8
$ws.Range('T130').Value2 = 2
$ws.Range('U130').Value2 = 3
$ws.Range('V130').Value2 = 1.95
$ws.Range('W130').Value2 = 1.85
$ws.Range('X130').Value2 = 0.363
$ws.Range('Y130').Value2 = -1
$ws.Range('Z130').Value2 = -1
$ws.Range('AA130').Value2 = 0.8
$ws.Range('AB130').Value2 = -1
$ws.Range('AC130').Value2 = 0
$ws.Range('AD130').Value2 = 0

# Row 131 <- original row 129 data (columns B:AD); column A (id) stays fixed
$ws.Range('B131').Value2 = 6782596
$ws.Range('C131').Value2 = 'Costa Rica Primera Division'
$ws.Range('D131').Value2 = 45255.95833333334
$ws.Range('E131').Value2 = 'Alajuelense'
$ws.Range('F131').Value2 = 'AD Guanacasteca'
$ws.Range('G131').Value2 = 3
$ws.Range('H131').Value2 = 4
$ws.Range('I131').Value2 = 0
$ws.Range('J131').Value2 = 3
$ws.Range('K131').Value2 = 'A'
$ws.Range('L131').Value2 = 1.363
$ws.Range('M131').Value2 = 4.75
$ws.Range('N131').Value2 = 8
$ws.Range('O131').Value2 = 1.444
$ws.Range('P131').Value2 = 4.333
$ws.Range('Q131').Value2 = 7
$ws.Range('R131').Value2 = -1.25
$ws.Range('S131').Value2 = 1.975
$ws.Range('T131').Value2 = 1.825
$ws.Range('U131').Value2 = 2.75
$ws.Range('V131').Value2 = 1.775
$ws.Range('W131').Value2 = 2.025
$ws.Range('X131').Value2 = -1
$ws.Range('Y131').Value2 = -1
$ws.Range('Z131').Value2 = 6
$ws.Range('AA131').Value2 = -1
$ws.Range('AB131').Value2 = 0.825
$ws.Range('AC131').Value2 = 0.7749999999999999
$ws.Range('AD131').Value2 = -1

# Row 224 <- original row 225 data (columns B:AD); column A (id) stays fixed
$ws.Range('B224').Value2 = 7623946
$ws.Range('C224').Value2 = 'Costa Rica Primera Division'
$ws.Range('D224').Value2 = 45388.83333333334
$ws.Range('E224').Value2 = 'Cartagines'
$ws.Range('F224').Value2 = 'Alajuelense'
$ws.Range('G224').Value2 = 0
$ws.Range('H224').Value2 = 0
$ws.Range('I224').Value2 = 0
$ws.Range('J224').Value2 = 0
$ws.Range('K224').Value2 = 'D'
$ws.Range('L224').Value2 = 3.4
$ws.Range('M224').Value2 = 3.4
$ws.Range('N224').Value2 = 1.95
$ws.Range('O224').Value2 = 3.8
$ws.Range('P224').Value2 = 3.6
$ws.Range('Q224').Value2 = 1.75
$ws.Range('R224').Value2 = 0.75
$ws.Range('S224').Value2 = 1.8
$ws.Range('T224').Value2 = 2
$ws.Range('U224').Value2 = 2.75
$ws.Range('V224').Value2 = 1.925
$ws.Range('W224').Value2 = 1.875
$ws.Range('X224').Value2 = -1
$ws.Range('Y224').Value2 = 2.6
$ws.Range('Z224').Value2 = -1
$ws.Range('AA224').Value2 = 0.8
$ws.Range('AB224').Value2 = -1
$ws.Range('AC224').Value2 = -1
$ws.Range('AD224').Value2 = 0.875

# Row 225 <- original row 224 data (columns B:AD); column A (id) stays fixed
$ws.Range('B225').Value2 = 7623944
$ws.Range('C225').Value2 = 'Costa Rica Primera Division'
$ws.Range('D225').Value2 = 45388.83333333334
$ws.Range('E225').Value2 = 'Santos de Gupiles'
$ws.Range('F225').Value2 = 'Municipal Liberia'
$ws.Range('G225').Value2 = 3
$ws.Range('H225').Value2 = 0
$ws.Range('I225').Value2 = 1
$ws.Range('J225').Value2 = 0
$ws.Range('K225').Value2 = 'H'
$ws.Range('L225').Value2 = 2.9
$ws.Range('M225').Value2 = 3.25
$ws.Range('N225').Value2 = 2.375
$ws.Range('O225').Value2 = 3.4
$ws.Range('P225').Value2 = 3.2
$ws.Range('Q225').Value2 = 2.15
$ws.Range('R225').Value2 = 0.25
$ws.Range('S225').Value2 = 1.95
$ws.Range('T225').Value2 = 1.85
$ws.Range('U225').Value2 = 2.5
$ws.Range('V225').Value2 = 2
$ws.Range('W225').Value2 = 1.8
$ws.Range('X225').Value2 = 2.4
$ws.Range('Y225').Value2 = -1
$ws.Range('Z225').Value2 = -1
$ws.Range('AA225').Value2 = 0.95
$ws.Range('AB225').Value2 = -1
$ws.Range('AC225').Value2 = 1
$ws.Range('AD225').Value2 = -1

# Row 267 <- original row 268 data (columns B:AD); column A (id) stays fixed
$ws.Range('B267').Value2 = 8162891
$ws.Range('C267').Value2 = 'Costa Rica Primera Division'
$ws.Range('D267').Value2 = 45424.75
$ws.Range('E267').Value2 = 'Deportivo Saprissa'
$ws.Range('F267').Value2 = 'Santos de Gupiles'
$ws.Range('G267').Value2 = 3
$ws.Range('H267').Value2 = 1
$ws.Range('I267').Value2 = 2
$ws.Range('J267').Value2 = 0
$ws.Range('K267').Value2 = 'H'
$ws.Range('L267').Value2 = 1.166
$ws.Range('M267').Value2 = 6.5
$ws.Range('N267').Value2 = 13
$ws.Range('O267').Value2 = 1.125
$ws.Range('P267').Value2 = 8
$ws.Range('Q267').Value2 = 15
$ws.Range('R267').Value2 = -2.25
$ws.Range('S267').Value2 = 1.825
$ws.Range('T267').Value2 = 1.975
$ws.Range('U267').Value2 = 3.5
$ws.Range('V267').Value2 = 1.975
$ws.Range('W267').Value2 = 1.825
$ws.Range('X267').Value2 = 0.125
$ws.Range('Y267').Value2 = -1
$ws.Range('Z267').Value2 = -1
$ws.Range('AA267').Value2 = -0.5
$ws.Range('AB267').Value2 = 0.4875
$ws.Range('AC267').Value2 = 0.9750000000000001
$ws.Range('AD267').Value2 = -1

# Row 268 <- original row 267 data (columns B:AD); column A (id) stays fixed
$ws.Range('B268').Value2 = 8203655
$ws.Range('C268').Value2 = 'Costa Rica Primera Division'
$ws.Range('D268').Value2 = 45424.75
$ws.Range('E268').Value2 = 'Municipal Perez Zeledon'
$ws.Range('F268').Value2 = 'Municipal Liberia'
$ws.Range('G268').Value2 = 0
$ws.Range('H268').Value2 = 3
$ws.Range('I268').Value2 = 0
$ws.Range('J268').Value2 = 1
$ws.Range('K268').Value2 = 'A'
$ws.Range('L268').Value2 = 3.3
$ws.Range('M268').Value2 = 3.5
$ws.Range('N268').Value2 = 2
$ws.Range('O268').Value2 = 2.9
$ws.Range('P268').Value2 = 3.4
$ws.Range('Q268').Value2 = 2.2
$ws.Range('R268').Value2 = 0.25
$ws.Range('S268').Value2 = 1.825
$ws.Range('T268').Value2 = 1.975
$ws.Range('U268').Value2 = 2.75
$ws.Range('V268').Value2 = 2
$ws.Range('W268').Value2 = 1.8
$ws.Range('X268').Value2 = -1
$ws.Range('Y268').Value2 = -1
$ws.Range('Z268').Value2 = 1.2
$ws.Range('AA268').Value2 = -1
$ws.Range('AB268').Value2 = 0.9750000000000001
$ws.Range('AC268').Value2 = 0.5
$ws.Range('AD268').Value2 = -0.5

# Row 269 <- original row 271 data (columns B:AD); column A (id) stays fixed
$ws.Range('B269').Value2 = 8162895
$ws.Range('C269').Value2 = 'Costa Rica Primera Division'
$ws.Range('D269').Value2 = 45424.75
$ws.Range('E269').Value2 = 'Sporting San Jose'
$ws.Range('F269').Value2 = 'Herediano'
$ws.Range('G269').Value2 = 1
$ws.Range('H269').Value2 = 1
$ws.Range('I269').Value2 = 1
$ws.Range('J269').Value2 = 0
$ws.Range('K269').Value2 = 'D'
$ws.Range('L269').Value2 = 3.6
$ws.Range('M269').Value2 = 3.5
$ws.Range('N269').Value2 = 1.833
$ws.Range('O269').Value2 = 4.5
$ws.Range('P269').Value2 = 3.8
$ws.Range('Q269').Value2 = 1.571
$ws.Range('R269').Value2 = 0.75
$ws.Range('S269').Value2 = 2.025
$ws.Range('T269').Value2 = 1.775
$ws.Range('U269').Value2 = 2.75
$ws.Range('V269').Value2 = 1.975
$ws.Range('W269').Value2 = 1.825
$ws.Range('X269').Value2 = -1
$ws.Range('Y269').Value2 = 2.8
$ws.Range('Z269').Value2 = -1
$ws.Range('AA269').Value2 = 1.025
$ws.Range('AB269').Value2 = -1
$ws.Range('AC269').Value2 = -1
$ws.Range('AD269').Value2 = 0.825

# Row 270 <- original row 269 data (columns B:AD); column A (id) stays fixed
$ws.Range('B270').Value2 = 8162892
$ws.Range('C270').Value2 = 'Costa Rica Primera Division'
$ws.Range('D270').Value2 = 45424.75
$ws.Range('E270').Value2 = 'Alajuelense'
$ws.Range('F270').Value2 = 'AD Guanacasteca'
$ws.Range('G270').Value2 = 5
$ws.Range('H270').Value2 = 0
$ws.Range('I270').Value2 = 2
$ws.Range('J270').Value2 = 0
$ws.Range('K270').Value2 = 'H'
$ws.Range('L270').Value2 = 1.25
$ws.Range('M270').Value2 = 5
$ws.Range('N270').Value2 = 10
$ws.Range('O270').Value2 = 1.3
$ws.Range('P270').Value2 = 4.75
$ws.Range('Q270').Value2 = 8
$ws.Range('R270').Value2 = -1.5
$ws.Range('S270').Value2 = 1.9
$ws.Range('T270').Value2 = 1.9
$ws.Range('U270').Value2 = 3
$ws.Range('V270').Value2 = 1.9
$ws.Range('W270').Value2 = 1.9
$ws.Range('X270').Value2 = 0.3
$ws.Range('Y270').Value2 = -1
$ws.Range('Z270').Value2 = -1
$ws.Range('AA270').Value2 = 0.8999999999999999
$ws.Range('AB270').Value2 = -1
$ws.Range('AC270').Value2 = 0.8999999999999999
$ws.Range('AD270').Value2 = -1

# Row 271 <- original row 270 data (columns B:AD); column A (id) stays fixed
$ws.Range('B271').Value2 = 8162893
$ws.Range('C271').Value2 = 'Costa Rica Primera Division'
$ws.Range('D271').Value2 = 45424.75
$ws.Range('E271').Value2 = 'AD Grecia'
$ws.Range('F271').Value2 = 'AD San Carlos'
$ws.Range('G271').Value2 = 2
$ws.Range('H271').Value2 = 2
$ws.Range('I271').Value2 = 0
$ws.Range('J271').Value2 = 1
$ws.Range('K271').Value2 = 'D'
$ws.Range('L271').Value2 = 5
$ws.Range('M271').Value2 = 4
$ws.Range('N271').Value2 = 1.533
$ws.Range('O271').Value2 = 4.2
$ws.Range('P271').Value2 = 4.2
$ws.Range('Q271').Value2 = 1.6
$ws.Range('R271').Value2 = 1
$ws.Range('S271').Value2 = 1.775
$ws.Range('T271').Value2 = 2.025
$ws.Range('U271').Value2 = 3
$ws.Range('V271').Value2 = 1.925
$ws.Range('W271').Value2 = 1.875
$ws.Range('X271').Value2 = -1
$ws.Range('Y271').Value2 = 3.2
$ws.Range('Z271').Value2 = -1
$ws.Range('AA271').Value2 = 0.7749999999999999
$ws.Range('AB271').Value2 = -1
$ws.Range('AC271').Value2 = 0.925
$ws.Range('AD271').Value2 = -1
